$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product")

# Widen column B to fit new product names
$ws.Columns.Item(2).ColumnWidth = 17.52

# Insert a new row above the existing KIVIK row so the new PAX entry becomes row 1
$ws.Rows.Item(1).Insert()

# Row 1: new PAX entry
$ws.Cells.Item(1, 1).Value = 19288134
$ws.Cells.Item(1, 2).Value = "PAX"

# Row 2 (former row 1) already has 99011429 / KIVIK, nothing to change there

# Row 3: new LINNMON/ALEX entry
$ws.Cells.Item(3, 1).Value = 69222616
$ws.Cells.Item(3, 2).Value = "LINNMON/ALEX"
